$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OrdenSalida")

# New values for rows 2-13 (columns A,B,G,H,S,T,AB)
$data = @(
    @{ Row=2;  S="20215631"; T=48  },
    @{ Row=3;  S="20202336"; T=12  },
    @{ Row=4;  S="20202335"; T=24  },
    @{ Row=5;  S="20215632"; T=228 },
    @{ Row=6;  S="20202347"; T=36  },
    @{ Row=7;  S="20287251"; T=24  },
    @{ Row=8;  S="20202309"; T=24  },
    @{ Row=9;  S="20202318"; T=24  },
    @{ Row=10; S="20202310"; T=48  },
    @{ Row=11; S="20215634"; T=36  },
    @{ Row=12; S="20202339"; T=144 },
    @{ Row=13; S="20287256"; T=48  }
)

foreach ($item in $data) {
    $r = $item.Row

    foreach ($col in @("A","B","G","H","S","AB")) {
        $ws.Range("$col$r").NumberFormat = "@"
    }

    $ws.Range("A$r").Value = "46988488"
    $ws.Range("B$r").Value = "46988488"
    $ws.Range("G$r").Value = "20250725"
    $ws.Range("H$r").Value = "20250805"
    $ws.Range("S$r").Value = $item.S
    $ws.Range("T$r").Value = $item.T
    $ws.Range("AB$r").Value = "5265"
}

# Remove now-unused rows 14-16 (data previously extended to row 16)
$ws.Range("A14:AB16").EntireRow.Delete()
